# Translate the ContosoLearn Competitor SWOT document from English to
# Spanish (es-ES), as produced by the "Juno: check in to OLPRODLOC."
# localization commit.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = $old
    $find.Replacement.ClearFormatting()
    $find.Replacement.Text = $new
    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null
}

function Replace-Split($old, $newA, $newB, $newC) {
    # Replace $old (a whole run's text) with $newA, then append $newB and
    # $newC as their own separate runs (mirrors how the source run gets
    # split into multiple <w:r> elements during localization).
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = $old
    $find.Replacement.ClearFormatting()
    $find.Replacement.Text = $newA
    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

    $r = $find.Parent
    $r.Collapse(0)
    $r.InsertAfter($newB)
    $r.Collapse(0)
    $r.InsertAfter($newC)
}

# --- Title ---
Replace-Text "ContosoLearn Competitor SWOT" "DAFO de competidores ContosoLearn"

# --- Bold section headings (each appears twice, once per vendor) ---
Replace-Text "Strengths:" "Fortalezas:"
Replace-Text "Weaknesses:" "Puntos débiles:"
Replace-Text "Opportunities:" "Oportunidades:"
Replace-Text "Threats:" "Amenazas:"

# --- Fabrikam Learning body text ---
Replace-Split (" Fabrikam Learning provides a comprehensive set of analytics and reporting tools. It ensures the continuous monitoring of teaching and learning activities, as well as pinpointing problematic areas that need to be addressed."), `
    " Fabrikam Learning proporciona un conjunto completo de herramientas de análisis e informes.", `
    " ", `
    "Garantiza la supervisión continua de las actividades de enseñanza y aprendizaje, así como la identificación de áreas problemáticas que deben abordarse."

Replace-Text " While Fabrikam Learning has robust reporting capabilities, it might be overwhelming for some users due to its comprehensive nature." " aunque Fabrikam Learning tiene funcionalidades de informes sólidas, puede ser abrumador para algunos usuarios debido a su naturaleza completa."

Replace-Split (" There is a growing demand for personalized learning experiences and data-driven recommendations. Fabrikam Learning can leverage its robust analytics and reporting tools to meet this demand."), `
    " hay una creciente demanda de experiencias de aprendizaje personalizadas y recomendaciones controladas por datos.", `
    " ", `
    "Fabrikam Learning puede aprovechar sus sólidas herramientas de análisis e informes para satisfacer esta demanda."

Replace-Split (" The eLearning market is highly competitive with many players offering similar features. Fabrikam Learning needs to continuously innovate to stay ahead."), `
    " el mercado de eLearning es altamente competitivo con muchos jugadores que ofrecen características similares.", `
    " ", `
    "Fabrikam Learning debe innovar continuamente para mantenerse a la vanguardia."

# --- AdatumLearn body text ---
Replace-Split (" AdatumLearn offers courses on business analysis techniques such as MOST and SWOT. This shows their commitment to providing valuable content to their users."), `
    " AdatumLearn ofrece cursos sobre técnicas de análisis de negocios como MOST y DAFO.", `
    " ", `
    "Esto muestra su compromiso de proporcionar contenido valioso a sus usuarios."

Replace-Split (" The information provided in their courses is a compilation of third-party generated information. This might not be as valuable as original content."), `
    " la información proporcionada en sus cursos es una compilación de información generada por terceros.", `
    " ", `
    "Esto podría no ser tan valioso como el contenido original."

Replace-Split (" AdatumLearn can create more original content to provide unique value to their users. They can also expand their course offerings to cover more topics."), `
    " AdatumLearn puede crear contenido más original para proporcionar un valor único a sus usuarios.", `
    " ", `
    "También puede ampliar sus ofertas de cursos para tratar más temas."

Replace-Split (' Like Fabrikam Learning, AdatumLearn also faces stiff competition in the eLearning market. They need to continuously improve their offerings to stay competitive."'), `
    " al igual que Fabrikam Learning, AdatumLearn también se enfrenta a una competencia rígida en el mercado de eLearning.", `
    " ", `
    ('Necesita mejorar continuamente su oferta para mantenerse competitivo".')

# --- Paragraph direction: explicitly mark every content paragraph as
# left-to-right (w:bidi val=0) like the localized document does. The
# trailing empty paragraph is left untouched. ---
for ($i = 1; $i -le 11; $i++) {
    $d.Paragraphs.Item($i).Format.ReadingOrder = 0
}
